$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2024-07-22 00:05:53"
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.001
$ws.Range("J3").Value = 0.05
$ws.Range("K3").Value = 0.003
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = 500
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 6
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 1000
$ws.Range("R3").Value = 5
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = 30
$ws.Range("U3").Value = 0.9090909090909091
$ws.Range("V3").Value = "./Data/Crupier.xlsx"
$ws.Range("W3").Value = 233000
$ws.Range("X3").Value = "No es Simulación"
